# Update "想去人数" (interest count) figures across the four sheets to the
# latest scraped values (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# -- 展览 (Exhibitions) --
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 238
$ws.Range("F6").Value  = 80
$ws.Range("F7").Value  = 809
$ws.Range("F8").Value  = 458
$ws.Range("F14").Value = 396
$ws.Range("F15").Value = 6422
$ws.Range("F19").Value = 7379
$ws.Range("F22").Value = 3342
$ws.Range("F23").Value = 772
$ws.Range("F24").Value = 850
$ws.Range("F25").Value = 4502
$ws.Range("F29").Value = 1411
$ws.Range("F30").Value = 137
$ws.Range("F33").Value = 1091
$ws.Range("F34").Value = 1522
$ws.Range("F35").Value = 2108

# -- 演出 (Performances) --
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 61
$ws.Range("F3").Value = 42
$ws.Range("F5").Value = 73

# -- 本地生活 (Local life) --
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 237
$ws.Range("F3").Value = 1196
$ws.Range("F4").Value = 67

# -- 全部类型 (All types) --
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 237
$ws.Range("F4").Value  = 1196
$ws.Range("F5").Value  = 67
$ws.Range("F8").Value  = 238
$ws.Range("F9").Value  = 80
$ws.Range("F10").Value = 809
$ws.Range("F11").Value = 458
$ws.Range("F14").Value = 61
$ws.Range("F18").Value = 396
$ws.Range("F19").Value = 6422
$ws.Range("F23").Value = 7379
$ws.Range("F26").Value = 3342
$ws.Range("F27").Value = 772
$ws.Range("F28").Value = 850
$ws.Range("F29").Value = 4502
$ws.Range("F31").Value = 42
$ws.Range("F34").Value = 1411
$ws.Range("F35").Value = 137
$ws.Range("F38").Value = 1091
$ws.Range("F39").Value = 1522
$ws.Range("F41").Value = 2108
$ws.Range("F43").Value = 73
